$wb = $excel.ActiveWorkbook

# ---- Sheet: Ring of Fire ----
$ws = $wb.Worksheets.Item("Ring of Fire")
$ws.Range("A4:N4").Copy()
$ws.Range("A21:N21").PasteSpecial(-4122)
$ws.Range("A21").Value = "GRAY"
$ws.Range("B21").Value = "MEG"
$ws.Range("C21").Value = "DRACO"
$ws.Range("D21").Value = "CARL"
$ws.Range("E21").Value = "8-BIT"
$ws.Range("F21").Value = "LUMI"
$ws.Range("G21").Value = "Equipo 2"
$ws.Range("H21").Value = "NXT|Rup"
$ws.Range("I21").Value = "NXT|amos"
$ws.Range("J21").Value = "NXT|Arthur"
$ws.Range("K21").Value = "IC|Mebius"
$ws.Range("L21").Value = "IC|RamaZR"
$ws.Range("M21").Value = "IC|Nob"
$ws.Range("N21").Value = "20250723T190839.000Z"

$ws.Range("A8:N8").Copy()
$ws.Range("A22:N22").PasteSpecial(-4122)
$ws.Range("A22").Value = "GRAY"
$ws.Range("B22").Value = "MEG"
$ws.Range("C22").Value = "DRACO"
$ws.Range("D22").Value = "CARL"
$ws.Range("E22").Value = "8-BIT"
$ws.Range("F22").Value = "LUMI"
$ws.Range("G22").Value = "Equipo 1"
$ws.Range("H22").Value = "NXT|Rup"
$ws.Range("I22").Value = "NXT|amos"
$ws.Range("J22").Value = "NXT|Arthur"
$ws.Range("K22").Value = "IC|Mebius"
$ws.Range("L22").Value = "IC|RamaZR"
$ws.Range("M22").Value = "IC|Nob"
$ws.Range("N22").Value = "20250723T190617.000Z"

# ---- Sheet: New Horizons ----
$ws = $wb.Worksheets.Item("New Horizons")
$ws.Range("A8:N8").Copy()
$ws.Range("A25:N25").PasteSpecial(-4122)
$ws.Range("A25").Value = "LUMI"
$ws.Range("B25").Value = "HANK"
$ws.Range("C25").Value = "MEG"
$ws.Range("D25").Value = "OLLIE"
$ws.Range("E25").Value = "MAX"
$ws.Range("F25").Value = "GENE"
$ws.Range("G25").Value = "Equipo 2"
$ws.Range("H25").Value = "FUT|GeRo"
$ws.Range("I25").Value = "FUT|Nowy297"
$ws.Range("J25").Value = "FUT|MeOw"
$ws.Range("K25").Value = "Enraged 💔"
$ws.Range("L25").Value = "SUP|Filippo神"
$ws.Range("M25").Value = "SUP|Tomzy"
$ws.Range("N25").Value = "20250723T192017.000Z"

$ws.Range("A8:N8").Copy()
$ws.Range("A26:N26").PasteSpecial(-4122)
$ws.Range("A26").Value = "LUMI"
$ws.Range("B26").Value = "HANK"
$ws.Range("C26").Value = "MEG"
$ws.Range("D26").Value = "OLLIE"
$ws.Range("E26").Value = "MAX"
$ws.Range("F26").Value = "GENE"
$ws.Range("G26").Value = "Equipo 2"
$ws.Range("H26").Value = "FUT|GeRo"
$ws.Range("I26").Value = "FUT|Nowy297"
$ws.Range("J26").Value = "FUT|MeOw"
$ws.Range("K26").Value = "Enraged 💔"
$ws.Range("L26").Value = "SUP|Filippo神"
$ws.Range("M26").Value = "SUP|Tomzy"
$ws.Range("N26").Value = "20250723T191721.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A27:N27").PasteSpecial(-4122)
$ws.Range("A27").Value = "LUMI"
$ws.Range("B27").Value = "HANK"
$ws.Range("C27").Value = "MEG"
$ws.Range("D27").Value = "OLLIE"
$ws.Range("E27").Value = "MAX"
$ws.Range("F27").Value = "GENE"
$ws.Range("G27").Value = "Equipo 1"
$ws.Range("H27").Value = "FUT|GeRo"
$ws.Range("I27").Value = "FUT|Nowy297"
$ws.Range("J27").Value = "FUT|MeOw"
$ws.Range("K27").Value = "Enraged 💔"
$ws.Range("L27").Value = "SUP|Filippo神"
$ws.Range("M27").Value = "SUP|Tomzy"
$ws.Range("N27").Value = "20250723T191316.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A28:N28").PasteSpecial(-4122)
$ws.Range("A28").Value = "BUSTER"
$ws.Range("B28").Value = "GRIFF"
$ws.Range("C28").Value = "BONNIE"
$ws.Range("D28").Value = "FRANK"
$ws.Range("E28").Value = "POCO"
$ws.Range("F28").Value = "GENE"
$ws.Range("G28").Value = "Equipo 1"
$ws.Range("H28").Value = "FUT|Nowy297"
$ws.Range("I28").Value = "FUT|MeOw"
$ws.Range("J28").Value = "FUT|GeRo"
$ws.Range("K28").Value = "Enraged 💔"
$ws.Range("L28").Value = "SUP|Tomzy"
$ws.Range("M28").Value = "SUP|Filippo神"
$ws.Range("N28").Value = "20250723T190632.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A29:N29").PasteSpecial(-4122)
$ws.Range("A29").Value = "KAZE"
$ws.Range("B29").Value = "BROCK"
$ws.Range("C29").Value = "BELLE"
$ws.Range("D29").Value = "HANK"
$ws.Range("E29").Value = "MEEPLE"
$ws.Range("F29").Value = "BYRON"
$ws.Range("G29").Value = "Equipo 1"
$ws.Range("H29").Value = "NOVO|Marco"
$ws.Range("I29").Value = "NOVO|Subeme"
$ws.Range("J29").Value = "NOVO|Biso"
$ws.Range("K29").Value = "TH|LeNain"
$ws.Range("L29").Value = "TH|Zhar"
$ws.Range("M29").Value = "TH|iKaoss"
$ws.Range("N29").Value = "20250723T193329.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A30:N30").PasteSpecial(-4122)
$ws.Range("A30").Value = "KAZE"
$ws.Range("B30").Value = "BROCK"
$ws.Range("C30").Value = "BELLE"
$ws.Range("D30").Value = "HANK"
$ws.Range("E30").Value = "MEEPLE"
$ws.Range("F30").Value = "BYRON"
$ws.Range("G30").Value = "Equipo 1"
$ws.Range("H30").Value = "NOVO|Marco"
$ws.Range("I30").Value = "NOVO|Subeme"
$ws.Range("J30").Value = "NOVO|Biso"
$ws.Range("K30").Value = "TH|LeNain"
$ws.Range("L30").Value = "TH|Zhar"
$ws.Range("M30").Value = "TH|iKaoss"
$ws.Range("N30").Value = "20250723T193122.000Z"

$ws.Range("A8:N8").Copy()
$ws.Range("A31:N31").PasteSpecial(-4122)
$ws.Range("A31").Value = "KAZE"
$ws.Range("B31").Value = "BROCK"
$ws.Range("C31").Value = "BELLE"
$ws.Range("D31").Value = "HANK"
$ws.Range("E31").Value = "MEEPLE"
$ws.Range("F31").Value = "BYRON"
$ws.Range("G31").Value = "Equipo 2"
$ws.Range("H31").Value = "NOVO|Marco"
$ws.Range("I31").Value = "NOVO|Subeme"
$ws.Range("J31").Value = "NOVO|Biso"
$ws.Range("K31").Value = "TH|LeNain"
$ws.Range("L31").Value = "TH|Zhar"
$ws.Range("M31").Value = "TH|iKaoss"
$ws.Range("N31").Value = "20250723T192757.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A32:N32").PasteSpecial(-4122)
$ws.Range("A32").Value = "CHARLIE"
$ws.Range("B32").Value = "GUS"
$ws.Range("C32").Value = "GRAY"
$ws.Range("D32").Value = "BONNIE"
$ws.Range("E32").Value = "ANGELO"
$ws.Range("F32").Value = "SQUEAK"
$ws.Range("G32").Value = "Equipo 1"
$ws.Range("H32").Value = "NOVO|Biso"
$ws.Range("I32").Value = "NOVO|Subeme"
$ws.Range("J32").Value = "NOVO|Marco"
$ws.Range("K32").Value = "TH|LeNain"
$ws.Range("L32").Value = "TH|iKaoss"
$ws.Range("M32").Value = "TH|Zhar"
$ws.Range("N32").Value = "20250723T192116.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A33:N33").PasteSpecial(-4122)
$ws.Range("A33").Value = "CHARLIE"
$ws.Range("B33").Value = "GUS"
$ws.Range("C33").Value = "GRAY"
$ws.Range("D33").Value = "BONNIE"
$ws.Range("E33").Value = "ANGELO"
$ws.Range("F33").Value = "SQUEAK"
$ws.Range("G33").Value = "Equipo 1"
$ws.Range("H33").Value = "NOVO|Biso"
$ws.Range("I33").Value = "NOVO|Subeme"
$ws.Range("J33").Value = "NOVO|Marco"
$ws.Range("K33").Value = "TH|LeNain"
$ws.Range("L33").Value = "TH|iKaoss"
$ws.Range("M33").Value = "TH|Zhar"
$ws.Range("N33").Value = "20250723T191844.000Z"

# ---- Sheet: Hot Potato ----
$ws = $wb.Worksheets.Item("Hot Potato")
$ws.Range("A6:N6").Copy()
$ws.Range("A40:N40").PasteSpecial(-4122)
$ws.Range("A40").Value = "CORDELIUS"
$ws.Range("B40").Value = "KAZE"
$ws.Range("C40").Value = "BEA"
$ws.Range("D40").Value = "RUFFS"
$ws.Range("E40").Value = "BERRY"
$ws.Range("F40").Value = "CROW"
$ws.Range("G40").Value = "Equipo 1"
$ws.Range("H40").Value = "NOVO|Subeme"
$ws.Range("I40").Value = "NOVO|Marco"
$ws.Range("J40").Value = "NOVO|Biso"
$ws.Range("K40").Value = "TH|LeNain"
$ws.Range("L40").Value = "TH|iKaoss"
$ws.Range("M40").Value = "TH|Zhar"
$ws.Range("N40").Value = "20250723T191034.000Z"

$ws.Range("A6:N6").Copy()
$ws.Range("A41:N41").PasteSpecial(-4122)
$ws.Range("A41").Value = "CORDELIUS"
$ws.Range("B41").Value = "KAZE"
$ws.Range("C41").Value = "BEA"
$ws.Range("D41").Value = "RUFFS"
$ws.Range("E41").Value = "BERRY"
$ws.Range("F41").Value = "CROW"
$ws.Range("G41").Value = "Equipo 1"
$ws.Range("H41").Value = "NOVO|Subeme"
$ws.Range("I41").Value = "NOVO|Marco"
$ws.Range("J41").Value = "NOVO|Biso"
$ws.Range("K41").Value = "TH|LeNain"
$ws.Range("L41").Value = "TH|iKaoss"
$ws.Range("M41").Value = "TH|Zhar"
$ws.Range("N41").Value = "20250723T190644.000Z"

# ---- Sheet: Dry Season ----
$ws = $wb.Worksheets.Item("Dry Season")
$ws.Range("A5:N5").Copy()
$ws.Range("A30:N30").PasteSpecial(-4122)
$ws.Range("A30").Value = "MORTIS"
$ws.Range("B30").Value = "MEEPLE"
$ws.Range("C30").Value = "GUS"
$ws.Range("D30").Value = "KAZE"
$ws.Range("E30").Value = "SQUEAK"
$ws.Range("F30").Value = "BELLE"
$ws.Range("G30").Value = "Equipo 2"
$ws.Range("H30").Value = "IC|Mebius"
$ws.Range("I30").Value = "IC|Nob"
$ws.Range("J30").Value = "IC|RamaZR"
$ws.Range("K30").Value = "NXT|Arthur"
$ws.Range("L30").Value = "NXT|Rup"
$ws.Range("M30").Value = "NXT|amos"
$ws.Range("N30").Value = "20250723T193558.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A31:N31").PasteSpecial(-4122)
$ws.Range("A31").Value = "MORTIS"
$ws.Range("B31").Value = "GUS"
$ws.Range("C31").Value = "MEEPLE"
$ws.Range("D31").Value = "SQUEAK"
$ws.Range("E31").Value = "BELLE"
$ws.Range("F31").Value = "KAZE"
$ws.Range("G31").Value = "Equipo 1"
$ws.Range("H31").Value = "IC|Mebius"
$ws.Range("I31").Value = "IC|RamaZR"
$ws.Range("J31").Value = "IC|Nob"
$ws.Range("K31").Value = "NXT|Rup"
$ws.Range("L31").Value = "NXT|amos"
$ws.Range("M31").Value = "NXT|Arthur"
$ws.Range("N31").Value = "20250723T193327.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A32:N32").PasteSpecial(-4122)
$ws.Range("A32").Value = "GENE"
$ws.Range("B32").Value = "MR. P"
$ws.Range("C32").Value = "BONNIE"
$ws.Range("D32").Value = "PENNY"
$ws.Range("E32").Value = "JAE-YONG"
$ws.Range("F32").Value = "GUS"
$ws.Range("G32").Value = "Equipo 1"
$ws.Range("H32").Value = "NXT|amos"
$ws.Range("I32").Value = "NXT|Arthur"
$ws.Range("J32").Value = "NXT|Rup"
$ws.Range("K32").Value = "IC|Mebius"
$ws.Range("L32").Value = "IC|Nob"
$ws.Range("M32").Value = "IC|RamaZR"
$ws.Range("N32").Value = "20250723T192116.000Z"

$ws.Range("A4:N4").Copy()
$ws.Range("A33:N33").PasteSpecial(-4122)
$ws.Range("A33").Value = "GENE"
$ws.Range("B33").Value = "MR. P"
$ws.Range("C33").Value = "BONNIE"
$ws.Range("D33").Value = "PENNY"
$ws.Range("E33").Value = "JAE-YONG"
$ws.Range("F33").Value = "GUS"
$ws.Range("G33").Value = "Equipo 1"
$ws.Range("H33").Value = "NXT|amos"
$ws.Range("I33").Value = "NXT|Arthur"
$ws.Range("J33").Value = "NXT|Rup"
$ws.Range("K33").Value = "IC|Mebius"
$ws.Range("L33").Value = "IC|Nob"
$ws.Range("M33").Value = "IC|RamaZR"
$ws.Range("N33").Value = "20250723T191857.000Z"

$ws.Range("A5:N5").Copy()
$ws.Range("A34:N34").PasteSpecial(-4122)
$ws.Range("A34").Value = "GENE"
$ws.Range("B34").Value = "MR. P"
$ws.Range("C34").Value = "BONNIE"
$ws.Range("D34").Value = "PENNY"
$ws.Range("E34").Value = "JAE-YONG"
$ws.Range("F34").Value = "GUS"
$ws.Range("G34").Value = "Equipo 2"
$ws.Range("H34").Value = "NXT|amos"
$ws.Range("I34").Value = "NXT|Arthur"
$ws.Range("J34").Value = "NXT|Rup"
$ws.Range("K34").Value = "IC|Mebius"
$ws.Range("L34").Value = "IC|Nob"
$ws.Range("M34").Value = "IC|RamaZR"
$ws.Range("N34").Value = "20250723T191637.000Z"

